# Updating the models with january production data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(6378,6260,6235,6153,6143,6056,6069,6008,6033,5985,6030,5993,6012,5999,6039,6017,6038,6008,6100,6105,6208,6205,6297,6418,6605,6745,6845,7008,7211,7400,7481,7647)
$bValues = @(46072.95833333334,46072.96875,46072.97916666666,46072.98958333334,46073,46073.01041666666,46073.02083333334,46073.03125,46073.04166666666,46073.05208333334,46073.0625,46073.07291666666,46073.08333333334,46073.09375,46073.10416666666,46073.11458333334,46073.125,46073.13541666666,46073.14583333334,46073.15625,46073.16666666666,46073.17708333334,46073.1875,46073.19791666666,46073.20833333334,46073.21875,46073.22916666666,46073.23958333334,46073.25,46073.26041666666,46073.27083333334,46073.28125)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

$wb.Save()
